$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Row 8 (the long "Agregar margen del 30%..." entry) becomes hidden.
$ws.Rows(8).Hidden = $true

# NOTE on ordering: new shared strings are appended to sharedStrings.xml
# in the order the values are first written, so we write "Abono de
# Amelia" before the longer "Poner NOTA PAGADA..." text to match the
# target's shared-string index order (21, then 22).

# 2) Insert 2 new rows right after the "Abono Martin" row (row 24) to
#    hold the new "Abono de Amelia" payment plus an extra blank spacer
#    row, pushing the two trailing blank rows further down. This runs
#    before the insert below so the row numbers here are the original
#    (pre-shift) ones; the later insert shifts these rows down to 28/29.
$ws.Rows("25:26").Insert()

$ws.Range("B25").Value = "Abono de Amelia"
$ws.Range("C25").Value = 5000

# 3) Insert 3 new rows before the spacer row (old row 13), pushing the
#    spacer + everything below it down by 3. Excel copies formatting
#    from the row above (row 12) into the new rows automatically.
$ws.Rows("13:15").Insert()

# Row 13 becomes a new content row describing further system changes,
# amount 1800.
$ws.Range("B13").Value = "Poner NOTA PAGADA cuando pagan una nota, Agregar el cambio en el ticket de Pagos, Agregar sonidos a la verificación de inventario físico, Agregar sección de cambio de contraseña y de código de autorización"
$ws.Range("C13").Value = 1800
$ws.Rows(13).RowHeight = 42.75

# Rows 14 and 15 stay blank (already carry the right style from the insert).

# 4) Move the selection to match the saved workbook state.
$ws.Range("G12").Select()
